# Daily attendance processing - 2026-01-22 21:38:18
# Swap the order of names in the "Recorded By" column (G) so that entries
# that read "System, dnasr281@gmail.com" become "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Text -eq $oldText) {
        $cell.Value = $newText
    }
}
